# Apply the "Catch orbits now that exceed large values over 10^25" update.
$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Probe I") ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("E2").Value = 10
$ws1.Range("F2").Value = -1
$ws1.Range("G2").Value = 1
$ws1.Range("H2").Value = 0

$ws1.Range("E3").Value = 25
$ws1.Range("F3").Value = 2
$ws1.Range("G3").Value = 2
$ws1.Range("H3").Value = 0

$ws1.Range("H5").Select()

# --- Sheet2 ("Probe II") ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("D2").Value = 17000

$ws2.Range("E2").Value = 10
$ws2.Range("F2").Value = 1
$ws2.Range("G2").Value = -1
$ws2.Range("H2").Value = 0

$ws2.Range("E3").Value = 25
$ws2.Range("F3").Value = 1
$ws2.Range("G3").Value = 1
$ws2.Range("H3").Value = 0

$ws2.Activate()
$ws2.Range("L12").Select()
